# Edit script: adds a "metadata" sheet after "data", updates F-column
# timestamps on "data" to reflect the re-run query time.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update the time_taken (column F) values on the "data" sheet -------
$timeUpdates = @(
    @{ Row = 2; Val = "2021-10-05 14:20:47.983545" },
    @{ Row = 3; Val = "2021-10-05 14:20:47.983553" },
    @{ Row = 4; Val = "2021-10-05 14:20:47.983556" },
    @{ Row = 5; Val = "2021-10-05 14:20:47.983559" },
    @{ Row = 6; Val = "2021-10-05 14:20:47.983561" },
    @{ Row = 7; Val = "2021-10-05 14:20:47.983564" },
    @{ Row = 8; Val = "2021-10-05 14:20:47.983567" },
    @{ Row = 9; Val = "2021-10-05 14:20:47.983569" },
    @{ Row = 10; Val = "2021-10-05 14:20:47.983572" },
    @{ Row = 11; Val = "2021-10-05 14:20:47.983574" },
    @{ Row = 12; Val = "2021-10-05 14:20:47.983577" },
    @{ Row = 13; Val = "2021-10-05 14:20:47.983579" },
    @{ Row = 14; Val = "2021-10-05 14:20:47.983582" },
    @{ Row = 15; Val = "2021-10-05 14:20:47.983584" },
    @{ Row = 16; Val = "2021-10-05 14:20:47.983587" },
    @{ Row = 17; Val = "2021-10-05 14:20:47.983589" },
    @{ Row = 18; Val = "2021-10-05 14:20:47.983592" },
    @{ Row = 19; Val = "2021-10-05 14:20:47.983595" },
    @{ Row = 20; Val = "2021-10-05 14:20:47.983597" },
    @{ Row = 21; Val = "2021-10-05 14:20:47.983600" },
    @{ Row = 22; Val = "2021-10-05 14:20:47.983602" },
    @{ Row = 23; Val = "2021-10-05 14:20:47.983605" },
    @{ Row = 24; Val = "2021-10-05 14:20:47.983607" },
    @{ Row = 25; Val = "2021-10-05 14:20:47.983610" },
    @{ Row = 26; Val = "2021-10-05 14:20:47.983613" },
    @{ Row = 27; Val = "2021-10-05 14:20:47.983615" },
    @{ Row = 28; Val = "2021-10-05 14:20:47.983618" },
    @{ Row = 29; Val = "2021-10-05 14:20:47.983620" },
    @{ Row = 30; Val = "2021-10-05 14:20:47.983623" },
    @{ Row = 31; Val = "2021-10-05 14:20:47.983625" },
    @{ Row = 32; Val = "2021-10-05 14:20:47.983628" },
    @{ Row = 33; Val = "2021-10-05 14:20:47.983630" },
    @{ Row = 34; Val = "2021-10-05 14:20:47.983633" },
    @{ Row = 35; Val = "2021-10-05 14:20:47.983636" },
    @{ Row = 36; Val = "2021-10-05 14:20:47.983638" },
    @{ Row = 37; Val = "2021-10-05 14:20:47.983641" },
    @{ Row = 38; Val = "2021-10-05 14:20:47.983643" },
    @{ Row = 39; Val = "2021-10-05 14:20:47.983646" },
    @{ Row = 40; Val = "2021-10-05 14:20:47.983648" },
    @{ Row = 41; Val = "2021-10-05 14:20:47.983651" },
    @{ Row = 42; Val = "2021-10-05 14:20:47.983654" },
    @{ Row = 43; Val = "2021-10-05 14:20:47.983656" },
    @{ Row = 44; Val = "2021-10-05 14:20:47.983659" },
    @{ Row = 45; Val = "2021-10-05 14:20:47.983661" },
    @{ Row = 46; Val = "2021-10-05 14:20:47.983664" },
    @{ Row = 47; Val = "2021-10-05 14:20:47.983666" },
    @{ Row = 48; Val = "2021-10-05 14:20:47.983669" },
    @{ Row = 49; Val = "2021-10-05 14:20:47.983671" },
    @{ Row = 50; Val = "2021-10-05 14:20:47.983673" },
    @{ Row = 51; Val = "2021-10-05 14:20:47.983676" },
    @{ Row = 52; Val = "2021-10-05 14:20:47.983678" },
    @{ Row = 53; Val = "2021-10-05 14:20:47.983681" },
    @{ Row = 54; Val = "2021-10-05 14:20:47.983684" },
    @{ Row = 55; Val = "2021-10-05 14:20:47.983686" },
    @{ Row = 56; Val = "2021-10-05 14:20:47.983689" },
    @{ Row = 57; Val = "2021-10-05 14:20:47.983691" },
    @{ Row = 58; Val = "2021-10-05 14:20:47.983694" },
    @{ Row = 59; Val = "2021-10-05 14:20:47.983696" },
    @{ Row = 60; Val = "2021-10-05 14:20:47.983699" },
    @{ Row = 61; Val = "2021-10-05 14:20:47.983701" },
    @{ Row = 62; Val = "2021-10-05 14:20:47.983704" },
    @{ Row = 63; Val = "2021-10-05 14:20:47.983706" },
    @{ Row = 64; Val = "2021-10-05 14:20:47.983709" },
    @{ Row = 65; Val = "2021-10-05 14:20:47.983711" },
    @{ Row = 66; Val = "2021-10-05 14:20:47.983714" },
    @{ Row = 67; Val = "2021-10-05 14:20:47.983717" },
    @{ Row = 68; Val = "2021-10-05 14:20:47.983720" },
    @{ Row = 69; Val = "2021-10-05 14:20:47.983722" },
    @{ Row = 70; Val = "2021-10-05 14:20:47.983725" },
    @{ Row = 71; Val = "2021-10-05 14:20:47.983727" },
    @{ Row = 72; Val = "2021-10-05 14:20:47.983729" },
    @{ Row = 73; Val = "2021-10-05 14:20:47.983732" },
    @{ Row = 74; Val = "2021-10-05 14:20:47.983735" },
    @{ Row = 75; Val = "2021-10-05 14:20:47.983737" },
    @{ Row = 76; Val = "2021-10-05 14:20:47.983740" },
    @{ Row = 77; Val = "2021-10-05 14:20:47.983743" },
    @{ Row = 78; Val = "2021-10-05 14:20:47.983747" },
    @{ Row = 79; Val = "2021-10-05 14:20:47.983750" },
    @{ Row = 80; Val = "2021-10-05 14:20:47.983753" },
    @{ Row = 81; Val = "2021-10-05 14:20:47.983755" },
    @{ Row = 82; Val = "2021-10-05 14:20:47.983758" },
    @{ Row = 83; Val = "2021-10-05 14:20:47.983760" },
    @{ Row = 84; Val = "2021-10-05 14:20:47.983763" },
    @{ Row = 85; Val = "2021-10-05 14:20:47.983765" },
    @{ Row = 86; Val = "2021-10-05 14:20:47.983767" },
    @{ Row = 87; Val = "2021-10-05 14:20:47.983770" },
    @{ Row = 88; Val = "2021-10-05 14:20:47.983772" },
    @{ Row = 89; Val = "2021-10-05 14:20:47.983775" },
    @{ Row = 90; Val = "2021-10-05 14:20:47.983777" },
    @{ Row = 91; Val = "2021-10-05 14:20:47.983780" },
    @{ Row = 92; Val = "2021-10-05 14:20:47.983782" },
    @{ Row = 93; Val = "2021-10-05 14:20:47.983784" },
    @{ Row = 94; Val = "2021-10-05 14:20:47.983788" },
    @{ Row = 95; Val = "2021-10-05 14:20:47.983791" },
    @{ Row = 96; Val = "2021-10-05 14:20:47.983794" },
    @{ Row = 97; Val = "2021-10-05 14:20:47.983796" },
    @{ Row = 98; Val = "2021-10-05 14:20:47.983799" },
    @{ Row = 99; Val = "2021-10-05 14:20:47.983801" },
    @{ Row = 100; Val = "2021-10-05 14:20:47.983804" },
    @{ Row = 101; Val = "2021-10-05 14:20:47.983806" },
    @{ Row = 102; Val = "2021-10-05 14:20:47.983808" },
    @{ Row = 103; Val = "2021-10-05 14:20:47.983811" },
    @{ Row = 104; Val = "2021-10-05 14:20:47.983813" },
    @{ Row = 105; Val = "2021-10-05 14:20:47.983816" }
)

foreach ($u in $timeUpdates) {
    $dataSheet.Cells.Item($u.Row, 6).Value = $u.Val
}

# --- 2. Add the new "metadata" sheet, placed after "data" -----------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Match the outline / page-setup conventions used on "data"
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1
$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

# Header row (row 1), columns B..G, bold + thin border + centered/top align
$headers = @(
    @{ Col = 2; Text = "data_name" },
    @{ Col = 3; Text = "data_id" },
    @{ Col = 4; Text = "data_version" },
    @{ Col = 5; Text = "data_version_created" },
    @{ Col = 6; Text = "panel_query_time" },
    @{ Col = 7; Text = "panel_get_request" }
)

foreach ($h in $headers) {
    $c = $metaSheet.Cells.Item(1, $h.Col)
    $c.Value = $h.Text
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

# Index cell A2 mirrors the styling used on the "data" sheet's A column
$idxCell = $metaSheet.Cells.Item(2, 1)
$idxCell.Value = 0
$idxCell.Font.Bold = $true
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160
$idxCell.Borders.LineStyle = 1

# Data row (row 2), columns B..G
$metaSheet.Cells.Item(2, 2).Value = "Hereditary spastic paraplegia - adult onset"
$metaSheet.Cells.Item(2, 3).Value = 567

# "1.72" must stay text (not be coerced to a number) -- a leading apostrophe
# forces Excel's quote-prefix text entry, same as typing '1.72 into a cell.
$metaSheet.Cells.Item(2, 4).Value = "'1.72"

$metaSheet.Cells.Item(2, 5).Value = "2021-09-01T11:02:23.913641Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:20:47.980177"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/567/?format=json"

# Leave "data" as the active sheet/tab, matching the original workbook view.
$dataSheet.Activate()

Write-Output "edit complete"
